$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Ativação:" date updated 01/01/2018 -> 01/01/2022 (row 8) ---
# Leading apostrophe keeps it text (matches original text-typed shared string,
# not an Excel date serial); re-paste formats from a same-column text cell so
# the cell style index (wrap/top-align/red-font) is preserved exactly.
$ws.Range("B8").Value = "'01/01/2022"
$ws.Range("C8").Value = "'01/01/2022"
$ws.Range("B3").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# --- New English "Objectives:" paragraph added (row 11, B/C previously empty) ---
$objectivesEn = "Provide the necessary knowledge on the fundamental aspects of Microbiology and Microbial Biochemistry and its importance in studies on Ecology of Microorganisms. Provide knowledge about the role and use of microorganisms in biological processes of interest to Environmental Engineering."
$ws.Range("B11").Value = $objectivesEn
$ws.Range("C11").Value = $objectivesEn
$ws.Range("B14").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- "Programa resumido:" text (row 14) - removed cellular structure/evolutive history clause ---
$programaResumido = "Diversidade metabólica; cultivo e crescimento microbiano; isolamento microbiano; ecossistemas microbianos; biorremediação e biodeterioração  microbiana; bioindicadores."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# --- "Short syllabus:" text (row 15) - removed cellular structure/evolutive history clause ---
$shortSyllabus = "Metabolic diversity; microbial culture and growth; microbial isolation; microbial ecosystems; microbial bioremediation and biodeterioration; bioindicators."
$ws.Range("B15").Value = $shortSyllabus
$ws.Range("C15").Value = $shortSyllabus

# --- "Programa:" text (row 16) - cellular structure clause replaced with organic molecules clause ---
$programa = "Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos,proteínas e ácidos nucleicos.–Diversidade metabólica: Micro-organismos autotróficos e heterotróficos; glicólise; fermentações; respiração; via das pentoses-fosfato; fotossíntese. –Cultivo e crescimento microbiano: Nutrição microbiana; meios de cultura; fatores ambientais; reprodução e crescimento; medidas e controle de crescimento microbiano. –Isolamento microbiano: Técnicas e meios de isolamento.–Ecossistemas microbianos: Diversidade microbiana e ciclos biogeoquímicos. –Biorremediação e biodeterioração microbiana: Lixiviação bacteriana de metais; bioacumulação e biotransformação microbiana de metais; biodegradação de materiais lignocelulósicos; biodegradação de hidrocarbonetos; biodeterioração de monumentos históricos. –Bioindicadores: Bioindicadores de qualidade de água, ar e solo."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- "Syllabus:" text (row 17) - cellular structure clause replaced with organic molecules clause ---
$syllabus = "Structure and function of the main organic molecules: carbohydrates, lipids,proteins and nucleic acids.Metabolic diversity: autotrophic and heterotrophic microorganisms, glycolysis; fermentations; respiration; pentose-phosphate pathway; photosynthesis. Microbial culture and growth: microbial nutrition; culture media; ambiental factors; reproduction and growth; measures and control of microbial growth.Microbial isolation: techniques and isolation media.Microbial ecosystems: microbial diversity and biogeochemical cycles.  Microbial biorremediation and biorremediation: bacterial leaching of metals; microbial bioaccumulation and biotransformation of metals; biodegradation of lignocellulosic materials; biodegradation of hydrocarbonets; biodeterioration of hystoric monuments. Bioindicators: bioindicators of the quality of water, air and soil."
$ws.Range("B17").Value = $syllabus
$ws.Range("C17").Value = $syllabus
